$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Name="139 Highett St Apartment Complex Richmond Outbreak"; Count=11},
    @{Name="3153 Sacred Heart Community St Kilda Tier 1B Outbreak"; Count=16},
    @{Name="3175 The Bays Aged Care Facility Hastings Outbreak"; Count=13},
    @{Name="3600 Belvedere Aged Care Noble Park Outbreak"; Count=37},
    @{Name="3612 BlueCross Glengowrie Outbreak"; Count=38},
    @{Name="3684 Homestyle Aged Care Langford Grange Cranbourne East Outbreak"; Count=33},
    @{Name="3980 Arcare Keysborough Aged Care Keysborough Outbreak"; Count=23},
    @{Name="4075 Ferndale Gardens Aged Care Services Bayswater North Outbreak"; Count=17},
    @{Name="4518 Regis Aged Care Fawkner Outbreak"; Count=15},
    @{Name="AW Window Transport Group Depot North Geelong Outbreak"; Count=10},
    @{Name="Allied Pinnacle Factory Altona North Outbreak"; Count=17},
    @{Name="Bespoke Childcare Dingley Village Outbreak"; Count=12},
    @{Name="Bread Solutions Braeside Outbreak"; Count=14},
    @{Name="CS Square Caroline Springs Outbreak"; Count=11},
    @{Name="Child's Play Early Learning Centre Tarneit Outbreak"; Count=11},
    @{Name="Community Kids Pascoe Vale Early Education Centre Pascoe Vale Outbreak"; Count=22},
    @{Name="Essential Caravans Somerton Outbreak"; Count=10},
    @{Name="Guardian Childcare Caulfield Outbreak"; Count=21},
    @{Name="Hello Fresh Warehouse Ravenhall Outbreak"; Count=17},
    @{Name="Inghams Enterprises Somerville Outbreak"; Count=24},
    @{Name="Lantmannen Unibake Australia Mordialloc Outbreak"; Count=20},
    @{Name="Launch Housing City Edge Crisis Accommodation South Melbourne Outbreak"; Count=11},
    @{Name="MacKillop Family Services Residential Facility Glenroy Outbreak"; Count=10},
    @{Name="Melbourne Custody Centre Melbourne Outbreak"; Count=14},
    @{Name="Melbourne Youth Justice Centre Parkville Outbreak"; Count=10},
    @{Name="Metro Tunnel Shine Domain Site Albert Road Construction Site South Melbourne Outbreak"; Count=10},
    @{Name="Monash Health Kingston Centre South 5 Outbreak"; Count=10},
    @{Name="Northern Health The Northern Hospital Epping Outbreak"; Count=20},
    @{Name="Oceania Meat Processors Laverton North Outbreak"; Count=17},
    @{Name="Public Housing 140 Brunswick Street Fitzroy Outbreak"; Count=14},
    @{Name="Public housing 33 Alfred Street North Melbourne Outbreak"; Count=11},
    @{Name="Shawlands Caravan Park Dandenong South Outbreak"; Count=17},
    @{Name="St Vincents Hospital Emergency Department Melbourne Outbreak"; Count=38},
    @{Name="Target Distribution Centre Truganina Outbreak"; Count=11},
    @{Name="The Robin Hood Inn Drouin West Outbreak"; Count=49},
    @{Name="The Toolshed Bar Private Event Noojee Outbreak"; Count=16},
    @{Name="Turosi Breakwater Outbreak"; Count=11},
    @{Name="Visy Recycling Springvale Outbreak"; Count=14},
    @{Name="Werribee Mercy Hospital Emergency Department Outbreak"; Count=34},
    @{Name="Western Health Footscray Hospital Emergency Department Outbreak"; Count=11},
    @{Name="Western Health Sunshine Hospital Emergency Department Outbreak"; Count=24}
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item.Name
    $ws.Cells.Item($row, 2).Value = $item.Count
    $row = $row + 1
}

Write-Output "Updated $($row - 2) rows"
